$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 9000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 9000
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H92").Value = 962.6667
$ws.Range("I92").Value = 823.4286
$ws.Range("K92").Value = 823.4286
$ws.Range("M92").Value = 424.5714

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 776520
$ws.Range("I2").Value = 982492.06
$ws.Range("K2").Value = 982492.06
$ws.Range("M2").Value = -982379.06
$ws.Range("H74").Value = 62508172
$ws.Range("I74").Value = 100006880
$ws.Range("K74").Value = 100006880
$ws.Range("M74").Value = -100006006
$ws.Range("H77").Value = 62508172
$ws.Range("I77").Value = 100006880
$ws.Range("K77").Value = 500034400
$ws.Range("M77").Value = -500030032
$ws.Range("H116").Value = 776520
$ws.Range("I116").Value = 982492.06
$ws.Range("K116").Value = 982492.06
$ws.Range("M116").Value = -980198.06
$ws.Range("H132").Value = 3457881.2
$ws.Range("I132").Value = 2176945
$ws.Range("J132").Value = 10004889
$ws.Range("K132").Value = 6530835
$ws.Range("L132").Value = 30014667
$ws.Range("M132").Value = -6528305
$ws.Range("N132").Value = -30019727

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 776520
$ws.Range("I3").Value = 982492.06
$ws.Range("K3").Value = 982492.06
$ws.Range("M3").Value = -982378.06

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 40999.668
$ws.Range("J28").Value = 40999.668
$ws.Range("L28").Value = 40999.668
$ws.Range("N28").Value = -41489.668
$ws.Range("H31").Value = 3721.426
$ws.Range("J31").Value = 10033.223
$ws.Range("L31").Value = 10033.223
$ws.Range("N31").Value = -10623.223
$ws.Range("H34").Value = 3721.426
$ws.Range("J34").Value = 10033.223
$ws.Range("L34").Value = 10033.223
$ws.Range("N34").Value = -10437.223
$ws.Range("H62").Value = 1999.5
$ws.Range("I62").Value = 1999.5
$ws.Range("K62").Value = 1999.5
$ws.Range("M62").Value = -1375.5
$ws.Range("H65").Value = 1999.5
$ws.Range("I65").Value = 1999.5
$ws.Range("K65").Value = 9997.5
$ws.Range("M65").Value = -6877.5
$ws.Range("H97").Value = 33332.332
$ws.Range("J97").Value = 33332.332
$ws.Range("L97").Value = 33332.332
$ws.Range("N97").Value = -35314.332
$ws.Range("H132").Value = 37038750
$ws.Range("I132").Value = 45456230
$ws.Range("K132").Value = 136368690
$ws.Range("M132").Value = -136366160

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 5075.25
$ws.Range("I50").Value = 5152.5
$ws.Range("K50").Value = 15457.5
$ws.Range("M50").Value = -14976.5
$ws.Range("H53").Value = 5075.25
$ws.Range("I53").Value = 5152.5
$ws.Range("K53").Value = 15457.5
$ws.Range("M53").Value = -14976.5
$ws.Range("H60").Value = 5472.933
$ws.Range("I60").Value = 178.8
$ws.Range("K60").Value = 536.4000000000001
$ws.Range("M60").Value = -285.4000000000001
$ws.Range("H131").Value = 2221.6365
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4004
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 4000
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 4004
$ws.Range("J73").Value = 4000
$ws.Range("L73").Value = 4000
$ws.Range("N73").Value = -5872
$ws.Range("H102").Value = 5628.263
$ws.Range("I102").Value = 3271.0625
$ws.Range("J102").Value = 18200
$ws.Range("K102").Value = 3271.0625
$ws.Range("L102").Value = 18200
$ws.Range("M102").Value = -1649.0625
$ws.Range("N102").Value = -21444
$ws.Range("H113").Value = 36007.39
$ws.Range("I113").Value = 43036.83
$ws.Range("J113").Value = 6885.4287
$ws.Range("K113").Value = 43036.83
$ws.Range("L113").Value = 6885.4287
$ws.Range("M113").Value = -40866.83
$ws.Range("N113").Value = -11225.4287

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 50051
$ws.Range("I56").Value = 50051
$ws.Range("K56").Value = 50051
$ws.Range("M56").Value = -49360
$ws.Range("H68").Value = 12507500
$ws.Range("I68").Value = 12507500
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 12507500
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -12506751
$ws.Range("N68").ClearContents()
$ws.Range("H70").Value = 40000
$ws.Range("I70").Value = 40000
$ws.Range("K70").Value = 40000
$ws.Range("M70").Value = -39730
$ws.Range("H71").Value = 12507500
$ws.Range("I71").Value = 12507500
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 62537500
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -62533756
$ws.Range("N71").ClearContents()
$ws.Range("H73").Value = 40000
$ws.Range("I73").Value = 40000
$ws.Range("K73").Value = 40000
$ws.Range("M73").Value = -39064
$ws.Range("H74").Value = 58613.25
$ws.Range("I74").Value = 44818
$ws.Range("K74").Value = 44818
$ws.Range("M74").Value = -43820
$ws.Range("H77").Value = 58613.25
$ws.Range("I77").Value = 44818
$ws.Range("K77").Value = 134454
$ws.Range("M77").Value = -129462
$ws.Range("H103").Value = 15900.25
$ws.Range("J103").Value = 15900.25
$ws.Range("L103").Value = 15900.25
$ws.Range("N103").Value = -18244.25
$ws.Range("H136").Value = 2356.2307
$ws.Range("I136").Value = 2300.5334
$ws.Range("J136").Value = 2432.182
$ws.Range("K136").Value = 6901.600199999999
$ws.Range("L136").Value = 7296.545999999999
$ws.Range("M136").Value = -4351.600199999999
$ws.Range("N136").Value = -12396.546

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1033
$ws.Range("J14").Value = 1200
$ws.Range("L14").Value = 1200
$ws.Range("N14").Value = -1536
$ws.Range("H27").Value = 48000
$ws.Range("J27").Value = 48000
$ws.Range("L27").Value = 48000
$ws.Range("N27").Value = -48138
$ws.Range("H46").Value = 54500
$ws.Range("J46").Value = 54500
$ws.Range("L46").Value = 54500
$ws.Range("N46").Value = -54962
$ws.Range("H62").Value = 5275.6924
$ws.Range("I62").Value = 3090
$ws.Range("K62").Value = 3090
$ws.Range("M62").Value = -2466
$ws.Range("H65").Value = 5275.6924
$ws.Range("I65").Value = 3090
$ws.Range("K65").Value = 15450
$ws.Range("M65").Value = -12330
$ws.Range("H75").Value = 79198.164
$ws.Range("I75").Value = 25063.334
$ws.Range("J75").Value = 133333
$ws.Range("K75").Value = 25063.334
$ws.Range("L75").Value = 133333
$ws.Range("M75").Value = -24127.334
$ws.Range("N75").Value = -135205
$ws.Range("H78").Value = 79198.164
$ws.Range("I78").Value = 25063.334
$ws.Range("J78").Value = 133333
$ws.Range("K78").Value = 75190.002
$ws.Range("L78").Value = 399999
$ws.Range("M78").Value = -70510.002
$ws.Range("N78").Value = -409359
$ws.Range("H96").Value = 1675
$ws.Range("I96").Value = 1350
$ws.Range("K96").Value = 1350
$ws.Range("M96").Value = 23
$ws.Range("H115").Value = 49249.332
$ws.Range("J115").Value = 49249.332
$ws.Range("L115").Value = 49249.332
$ws.Range("N115").Value = -52383.332
$ws.Range("H122").Value = 1594.6
$ws.Range("I122").Value = 1510.0834
$ws.Range("K122").Value = 4530.2502
$ws.Range("M122").Value = -2080.2502
$ws.Range("H126").Value = 1111.25
$ws.Range("I126").Value = 1125
$ws.Range("J126").Value = 1097.5
$ws.Range("K126").Value = 3375
$ws.Range("L126").Value = 3292.5
$ws.Range("M126").Value = -905
$ws.Range("N126").Value = -8232.5
$ws.Range("H132").Value = 31253126
$ws.Range("I132").Value = 50002836
$ws.Range("J132").Value = 3608.5
$ws.Range("K132").Value = 150008508
$ws.Range("L132").Value = 10825.5
$ws.Range("M132").Value = -150005978
$ws.Range("N132").Value = -15885.5
$ws.Range("H134").Value = 54500
$ws.Range("J134").Value = 54500
$ws.Range("L134").Value = 163500
$ws.Range("N134").Value = -168570
$ws.Range("H136").Value = 23811674
